$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Preserve the existing "AFDB" source citation; it moves down to rows 29-30 ---
$afdbName = $ws.Range("A23").Text
$afdbCite = $ws.Range("A24").Text

$ws.Range("A29").Value = $afdbName
$ws.Range("A29").Font.Bold = $true

$ws.Range("A30").Value = $afdbCite
$ws.Range("A30").Font.Italic = $true

# --- New MSME size-classification table (rows 20-24) ---

# Header row: column labels (bold, matching the sheet's other title rows)
$ws.Range("B20").Value = "Number of employees"
$ws.Range("C20").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D20").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B20:D20").Font.Bold = $true

# Micro row
$ws.Range("A21").Value = "Micro"
$ws.Range("B21").Value = 'Micro: 1-4<br/>Very Small EnTurnover (VSE): 5-9'
$ws.Range("C21").Value = ""
$ws.Range("D21").Value = ""

# Small row
$ws.Range("A22").Value = "Small"
$ws.Range("B22").Value = "10-100"
$ws.Range("C22").Value = ""
$ws.Range("D22").Value = ""

# Medium row (reuses/overwrites the old A23 "AFDB" cell with fresh, unstyled content)
$ws.Range("A23").Value = "Medium"
$ws.Range("A23").Font.Bold = $false
$ws.Range("B23").Value = "100-500"
$ws.Range("C23").Value = ""
$ws.Range("D23").Value = ""

# Large row (reuses/overwrites the old A24 citation cell with fresh, unstyled content)
$ws.Range("A24").Value = "Large"
$ws.Range("A24").Font.Italic = $false
$ws.Range("B24").Value = ">500"
$ws.Range("C24").Value = ""
$ws.Range("D24").Value = ""
